$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.543.41'
$ws.Range('E2').Value = '  +3.94%  '
$ws.Range('D3').Value = '2.266.34'
$ws.Range('E3').Value = '  +1.33%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.97%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.623'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '63.38'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.99%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.422'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0996'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +10.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '57.30'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.25%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '25.60'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +12.48%  '
$ws.Range('E13').Value = '  -0.51%  '
$ws.Range('D14').Value = '2.603.58'
$ws.Range('E14').Value = '  +1.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.57'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.05%  '
$ws.Range('E16').Value = '  +3.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.808'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('D18').Value = '2.250.69'
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('D19').Value = '43.521.91'
$ws.Range('E19').Value = '  +4.09%  '
$ws.Range('D20').Value = '0.0₃0967'
$ws.Range('E20').Value = '  +5.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.79'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.06'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.76%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '247.47'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.35%  '
$ws.Range('E24').Value = '  +0.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.48'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.26'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.78'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '171.46'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.92%  '
$ws.Range('E29').Value = '  -3.43%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.45'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.48%  '
$ws.Range('E31').Value = '  +1.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.78'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.90%  '
$ws.Range('E33').Value = '  -0.76%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0685'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.97%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.06'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.58%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.67'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.78'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.37%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.66'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.78%  '
$ws.Range('E39').Value = '  -4.11%  '
$ws.Range('E40').Value = '  +2.11%  '
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('E42').Value = '  -4.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.50'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +18.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0960'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.09%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '17.08'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.03%  '
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.20'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.94%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '1.468.42'
$ws.Range('E47').Value = '  -0.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '96.30'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.05%  '
$ws.Range('B49').Value = 'FTXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.34'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.45%  '
$ws.Range('B50').Value = 'TerraClassic'
$ws.Range('C50').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000208'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -14.98%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.34'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.18%  '
